# sample1.xlsx — add "No.Of Pieces / BAG WEIGHT / PACKET WEIGHT" columns,
# drop "VALUATION METHOD", move "SALES A/C" ahead of "SALES RETURN A/C",
# and add the item-master dropdown/range validations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural column changes -------------------------------------------

# Insert 3 new (blank) columns right before the old "VALUATION METHOD"
# column (F). They become the new F:H.
$ws.Range("F1:H1").EntireColumn.Insert()

# The old "VALUATION METHOD" column is now shifted to I — remove it.
$ws.Range("I1").EntireColumn.Delete()

# Populate the headers for the 3 newly inserted columns.
$ws.Range("F1").Value2 = "No.Of Pieces"
$ws.Range("G1").Value2 = "BAG WEIGHT"
$ws.Range("H1").Value2 = "PACKET WEIGHT"

# Move "SALES A/C" so it sits right before "SALES RETURN A/C".
$ws.Range("R1").Value2 = "SALES A/C"
$ws.Range("S1").Value2 = "SALES RETURN A/C"
$ws.Range("T1").Value2 = "EAN NO"
$ws.Range("U1").Value2 = "HSN/SAC"

# --- Column widths ---------------------------------------------------------

$ws.Columns("A").ColumnWidth = 14.833333333333334
$ws.Columns("B").ColumnWidth = 15.333333333333334
$ws.Columns("E").ColumnWidth = 10.666666666666666
$ws.Columns("F").ColumnWidth = 13.666666666666666
$ws.Columns("G").ColumnWidth = 19.666666666666668
$ws.Columns("H").ColumnWidth = 19.666666666666668
$ws.Columns("I").ColumnWidth = 26.833333333333332
$ws.Columns("J").ColumnWidth = 31.333333333333332
$ws.Columns("K").ColumnWidth = 22.833333333333332
$ws.Columns("M").ColumnWidth = 19.333333333333332
$ws.Columns("N").ColumnWidth = 15.333333333333334
$ws.Columns("O").ColumnWidth = 9.166666666666666
$ws.Columns("P").ColumnWidth = 17.166666666666668
$ws.Columns("Q").ColumnWidth = 12.333333333333334
$ws.Columns("R").ColumnWidth = 12.333333333333334
$ws.Columns("S").ColumnWidth = 16.666666666666668
$ws.Columns("T").ColumnWidth = 11.166666666666666

# --- Data validations --------------------------------------------------

# TYPE column dropdown
$ws.Range("E1:E1048576").Validation.Add(3, 1, 1, '"Consumed,Finished Goods,Packing Material,Raw Material,By Product,By Product Packets,Wastage"')

# No.Of Pieces: whole number range
$ws.Range("F2").Validation.Add(1, 1, 1, 1, 100000000)

# SOURCE column dropdown
$ws.Range("L1:L1048576").Validation.Add(3, 1, 1, '"Produced,Purchased,Produced or Purchased"')

# USAGE column dropdown
$ws.Range("M1:M1048576").Validation.Add(3, 1, 1, '"General Consumption,Sale,Rejected,Produced or Sale,Rejected or Sale,Produced or Rejected,Produced or Sale or Rejected"')

# BAG WEIGHT: decimal range
$ws.Range("G2").Validation.Add(2, 1, 1, 1, 1000000000000)

# PACKET WEIGHT: decimal range
$ws.Range("H2").Validation.Add(2, 1, 1, 1, 10000000000)

# --- Selection ---------------------------------------------------------

$ws.Range("B5").Select() | Out-Null
